# Refresh the "cryptos" price list (GitHub Actions style data refresh).
# Note: several Price cells (column D) look like plain decimals to Excel's
# auto-detection (e.g. "1.00", "565.72") and would silently be converted to
# numbers, losing the original text formatting used throughout the sheet.
# A leading apostrophe forces those specific assignments to stay text,
# matching the original column D values, which are all stored as strings.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.559.30"
$ws.Range("E2").Value = "  -2.10%  "
$ws.Range("D3").Value = "2.880.87"
$ws.Range("E3").Value = "  -2.28%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'565.72"
$ws.Range("E5").Value = "  -4.56%  "
$ws.Range("D6").Value = "'142.01"
$ws.Range("E6").Value = "  -3.60%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "'0.499"
$ws.Range("E8").Value = "  -1.19%  "
$ws.Range("D9").Value = "2.878.11"
$ws.Range("E9").Value = "  -2.30%  "
$ws.Range("D10").Value = "'6.86"
$ws.Range("E10").Value = "  -2.55%  "
$ws.Range("E11").Value = "  -2.69%  "
$ws.Range("D12").Value = "'0.427"
$ws.Range("E12").Value = "  -2.21%  "
$ws.Range("D13").Value = "'0.0000229"
$ws.Range("E13").Value = "  -1.94%  "
$ws.Range("D14").Value = "'31.52"
$ws.Range("E14").Value = "  -2.93%  "
$ws.Range("E15").Value = "  -0.28%  "
$ws.Range("D16").Value = "3.360.42"
$ws.Range("E16").Value = "  -2.22%  "
$ws.Range("D17").Value = "61.531.25"
$ws.Range("E17").Value = "  -2.13%  "
$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").Value = "'6.47"
$ws.Range("E18").Value = "  -3.06%  "
$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").Value = "2.868.37"
$ws.Range("E19").Value = "  -2.74%  "
$ws.Range("D20").Value = "'428.50"
$ws.Range("E20").Value = "  -2.13%  "
$ws.Range("D21").Value = "'12.97"
$ws.Range("E21").Value = "  -3.14%  "
$ws.Range("D22").Value = "'0.648"
$ws.Range("E22").Value = "  -2.33%  "
$ws.Range("D23").Value = "'6.77"
$ws.Range("E23").Value = "  -3.18%  "
$ws.Range("D24").Value = "'78.70"
$ws.Range("E24").Value = "  -2.40%  "
$ws.Range("D25").Value = "'11.78"
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("D27").Value = "'9.92"
$ws.Range("E27").Value = "  -11.87%  "
$ws.Range("E28").Value = "  -6.02%  "
$ws.Range("E29").Value = "  +6.02%  "
$ws.Range("E30").Value = "  -3.89%  "
$ws.Range("E31").Value = "  -4.98%  "
$ws.Range("D32").Value = "'2.01"
$ws.Range("E32").Value = "  -9.67%  "
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("D35").Value = "'25.29"
$ws.Range("E35").Value = "  -4.00%  "
$ws.Range("E36").Value = "  -3.75%  "
$ws.Range("D37").Value = "'5.33"
$ws.Range("E37").Value = "  -4.82%  "
$ws.Range("D38").Value = "'48.76"
$ws.Range("E38").Value = "  -1.65%  "
$ws.Range("D39").Value = "'2.79"
$ws.Range("E39").Value = "  -7.18%  "
$ws.Range("E40").Value = "  -6.23%  "
$ws.Range("D41").Value = "'8.12"
$ws.Range("E41").Value = "  -3.47%  "
$ws.Range("E42").Value = "  -4.31%  "
$ws.Range("D43").Value = "'39.10"
$ws.Range("E43").Value = "  -1.10%  "
$ws.Range("D44").Value = "'0.264"
$ws.Range("E44").Value = "  -5.34%  "
$ws.Range("D45").Value = "2.674.46"
$ws.Range("E45").Value = "  -0.55%  "
$ws.Range("D46").Value = "'132.30"
$ws.Range("E46").Value = "  -2.18%  "
$ws.Range("E47").Value = "  -1.28%  "
$ws.Range("D49").Value = "'341.34"
$ws.Range("E49").Value = "  -4.67%  "
$ws.Range("E50").Value = "  -2.09%  "
$ws.Range("D51").Value = "'21.30"
$ws.Range("E51").Value = "  -5.93%  "
